$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the "Play CyberCatz Online Slot Game for Free" heading.
# -----------------------------------------------------------------------
$metaRange = $d.Content
$metaFound = $metaRange.Find.Execute("Meta description", $true, $false, $false, `
                                      $false, $false, $true, 1, $false, "", 0)
if ($metaFound) {
    $metaPara = $metaRange.Paragraphs(1)
    $metaPara.Range.Delete()
}

# -----------------------------------------------------------------------
# 2) Locate the closing paragraph that holds the old feature-image prompt
#    ("Create an eye-catching feature image for ..."), insert a new bold
#    "Play CyberCatz Online Slot Game for Free" paragraph right before it,
#    and then turn the old paragraph's own text into the meta-description
#    sentence (keeping its italic formatting).
# -----------------------------------------------------------------------
$imgRange = $d.Content
$imgFound = $imgRange.Find.Execute("Create an eye-catching feature image", $true, `
                                    $false, $false, $false, $false, $true, 1, `
                                    $false, "", 0)
if ($imgFound) {
    $imgPara = $imgRange.Paragraphs(1)

    $headingText = "Play CyberCatz Online Slot Game for Free"
    $insertPos = $imgPara.Range.Start
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertAfter($headingText + [char]13)

    $headingRange = $d.Range($insertPos, $insertPos + $headingText.Length)
    $headingRange.Font.Bold = $true

    # The old paragraph (now shifted one paragraph mark later) gets new text.
    $metaSentenceRange = $d.Content
    $metaSentenceFound = $metaSentenceRange.Find.Execute("Create an eye-catching feature image", `
                                                          $true, $false, $false, $false, $false, `
                                                          $true, 1, $false, "", 0)
    if ($metaSentenceFound) {
        $metaSentencePara = $metaSentenceRange.Paragraphs(1)
        $metaSentencePara.Range.Text = "Join the adventure in CyberCatz, a futuristic intergalactic online slot game. Play for free and trigger up to 30 free spins with the Free Spins feature."
    }
}
